$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.262.52'
$ws.Range("E2").Value = '  +1.72%  '

$ws.Range("D3").Value = '3.785.98'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '668.22'
$ws.Range("E5").Value = '  +6.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.55'
$ws.Range("E6").Value = '  +1.19%  '

$ws.Range("D7").Value = '3.784.07'
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  +0.91%  '

$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("E11").Value = '  +1.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.05'
$ws.Range("E12").Value = '  +5.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000243'
$ws.Range("E13").Value = '  -1.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.58'
$ws.Range("E14").Value = '  -0.28%  '

$ws.Range("D15").Value = '4.418.26'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '3.790.03'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").Value = '70.136.43'
$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("E20").Value = '  +0.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.39'
$ws.Range("E21").Value = '  +18.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '473.15'
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.711'
$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.57'
$ws.Range("E24").Value = '  -0.67%  '

$ws.Range("E25").Value = '  -4.27%  '

$ws.Range("E26").Value = '  +1.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.33'
$ws.Range("E27").Value = '  +2.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.11'
$ws.Range("E28").Value = '  -2.30%  '

$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").Value = '3.934.12'
$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.84'
$ws.Range("E31").Value = '  +6.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.31'
$ws.Range("E32").Value = '  +2.74%  '

$ws.Range("E33").Value = '  +2.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.42'
$ws.Range("E34").Value = '  +2.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.178'
$ws.Range("E35").Value = '  +8.92%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.09'
$ws.Range("E37").Value = '  +1.12%  '

$ws.Range("D38").Value = '3.737.74'
$ws.Range("E38").Value = '  +0.31%  '

$ws.Range("E39").Value = '  +0.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.38'
$ws.Range("E40").Value = '  -2.01%  '

$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.31%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.963'
$ws.Range("E43").Value = '  -0.53%  '

$ws.Range("E44").Value = '  +10.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.56'
$ws.Range("E46").Value = '  +5.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.76'
$ws.Range("E47").Value = '  +4.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.00'
$ws.Range("E48").Value = '  +2.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.43'
$ws.Range("E49").Value = '  +4.80%  '

$ws.Range("E50").Value = '  +0.53%  '

$ws.Range("E51").Value = '  +1.00%  '
